# Insert a new weekly record for "Macroferia Regional de Talca - Coliflor".
#
# The sheet is a long append-only log of weekly price observations, one row
# per week, sorted with the most-recent week at the bottom in this workbook
# edit (a new observation is inserted right before the existing row for the
# "Primera" quality record at row 289, pushing every subsequent row down by
# one). The new row carries the same fixed descriptive fields as its
# neighbours plus the new week's date and prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 289; everything from the old row 289 down to the old
# row 403 shifts down to 290-404 (the engine also grows the used range /
# dimension to A1:R404 automatically).
$ws.Rows(289).Insert()

# Populate the newly inserted row 289 with this week's observation.
$ws.Range("A289").Value = 5
$ws.Range("B289").Value = "Macroferia Regional de Talca"
$ws.Range("C289").Value = "Maule"
$ws.Range("D289").Value = 45027
$ws.Range("E289").Value = 7
$ws.Range("F289").Value = 100112008
$ws.Range("G289").Value = "Coliflor"
$ws.Range("H289").Value = "Sin especificar"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 3000
$ws.Range("K289").Value = 1000
$ws.Range("L289").Value = 1000
$ws.Range("M289").Value = 1000
$ws.Range("N289").Value = "$/unidad"
$ws.Range("O289").Value = "Región del Maule"
$ws.Range("P289").Value = 1000
$ws.Range("Q289").Value = 1
$ws.Range("R289").Value = "Hortaliza"
